$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths for columns A and B (per diff: col1 width=20.0391, col2 width=13.7734)
$ws.Columns.Item(1).ColumnWidth = 19.2
$ws.Columns.Item(2).ColumnWidth = 12.95

# Copy formatting (styles/borders/number-format) from the last existing data row (494)
# down through the new rows (495:520), then fill in values + row heights.
$ws.Range("A494:F494").Copy()
$ws.Range("A495:F520").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Rows.Item(495).RowHeight = 13.55
$ws.Cells.Item(495, 1).Value = "2024-12-02 01:14:20"
$ws.Cells.Item(495, 2).Value = "021267.NC"
$ws.Cells.Item(495, 3).Value = 221
$ws.Cells.Item(495, 4).Value = 191
$ws.Cells.Item(495, 5).Value = 14
$ws.Cells.Item(495, 6).Value = 1

$ws.Rows.Item(496).RowHeight = 13.55
$ws.Cells.Item(496, 1).Value = "2024-12-03 08:57:56"
$ws.Cells.Item(496, 2).Value = "21455.NC"
$ws.Cells.Item(496, 3).Value = 13
$ws.Cells.Item(496, 4).Value = 8
$ws.Cells.Item(496, 5).Value = 20
$ws.Cells.Item(496, 6).Value = 1

$ws.Rows.Item(497).RowHeight = 13.55
$ws.Cells.Item(497, 1).Value = "2024-12-03 09:51:02"
$ws.Cells.Item(497, 2).Value = "021461.NC"
$ws.Cells.Item(497, 3).Value = 9
$ws.Cells.Item(497, 4).Value = 4
$ws.Cells.Item(497, 5).Value = 25
$ws.Cells.Item(497, 6).Value = 1

$ws.Rows.Item(498).RowHeight = 13.55
$ws.Cells.Item(498, 1).Value = "2024-12-03 10:22:09"
$ws.Cells.Item(498, 2).Value = "021499.NC"
$ws.Cells.Item(498, 3).Value = 20
$ws.Cells.Item(498, 4).Value = 127
$ws.Cells.Item(498, 5).Value = 20
$ws.Cells.Item(498, 6).Value = 1

$ws.Rows.Item(499).RowHeight = 13.55
$ws.Cells.Item(499, 1).Value = "2024-12-04 15:13:47"
$ws.Cells.Item(499, 2).Value = "021522.NC"
$ws.Cells.Item(499, 3).Value = 23
$ws.Cells.Item(499, 4).Value = 153
$ws.Cells.Item(499, 5).Value = 10
$ws.Cells.Item(499, 6).Value = 1

$ws.Rows.Item(500).RowHeight = 13.55
$ws.Cells.Item(500, 1).Value = "2024-12-04 16:38:03"
$ws.Cells.Item(500, 2).Value = "021547.NC"
$ws.Cells.Item(500, 3).Value = 9
$ws.Cells.Item(500, 4).Value = 248
$ws.Cells.Item(500, 5).Value = 12
$ws.Cells.Item(500, 6).Value = 1

$ws.Rows.Item(501).RowHeight = 13.55
$ws.Cells.Item(501, 1).Value = "2024-12-04 17:41:04"
$ws.Cells.Item(501, 2).Value = "020922.NC"
$ws.Cells.Item(501, 3).Value = 361
$ws.Cells.Item(501, 4).Value = 51
$ws.Cells.Item(501, 5).Value = 21
$ws.Cells.Item(501, 6).Value = 1

$ws.Rows.Item(502).RowHeight = 13.55
$ws.Cells.Item(502, 1).Value = "2024-12-05 10:52:40"
$ws.Cells.Item(502, 2).Value = "021665.NC"
$ws.Cells.Item(502, 3).Value = 3
$ws.Cells.Item(502, 4).Value = 9
$ws.Cells.Item(502, 5).Value = 12
$ws.Cells.Item(502, 6).Value = 1

$ws.Rows.Item(503).RowHeight = 13.55
$ws.Cells.Item(503, 1).Value = "2024-12-05 11:34:40"
$ws.Cells.Item(503, 2).Value = "021666.NC"
$ws.Cells.Item(503, 3).Value = 41
$ws.Cells.Item(503, 4).Value = 125
$ws.Cells.Item(503, 5).Value = 2
$ws.Cells.Item(503, 6).Value = 1

$ws.Rows.Item(504).RowHeight = 13.55
$ws.Cells.Item(504, 1).Value = "2024-12-05 12:30:07"
$ws.Cells.Item(504, 2).Value = "021667.NC"
$ws.Cells.Item(504, 3).Value = 44
$ws.Cells.Item(504, 4).Value = 125
$ws.Cells.Item(504, 5).Value = 2
$ws.Cells.Item(504, 6).Value = 1

$ws.Rows.Item(505).RowHeight = 13.55
$ws.Cells.Item(505, 1).Value = "2024-12-05 13:26:26"
$ws.Cells.Item(505, 2).Value = "021670.NC"
$ws.Cells.Item(505, 3).Value = 40
$ws.Cells.Item(505, 4).Value = 120
$ws.Cells.Item(505, 5).Value = 2
$ws.Cells.Item(505, 6).Value = 1

$ws.Rows.Item(506).RowHeight = 13.55
$ws.Cells.Item(506, 1).Value = "2024-12-05 14:11:46"
$ws.Cells.Item(506, 2).Value = "021668.NC"
$ws.Cells.Item(506, 3).Value = 20
$ws.Cells.Item(506, 4).Value = 88
$ws.Cells.Item(506, 5).Value = 2
$ws.Cells.Item(506, 6).Value = 1

$ws.Rows.Item(507).RowHeight = 13.55
$ws.Cells.Item(507, 1).Value = "2024-12-05 14:44:26"
$ws.Cells.Item(507, 2).Value = "021669.NC"
$ws.Cells.Item(507, 3).Value = 11
$ws.Cells.Item(507, 4).Value = 36
$ws.Cells.Item(507, 5).Value = 2
$ws.Cells.Item(507, 6).Value = 1

$ws.Rows.Item(508).RowHeight = 13.55
$ws.Cells.Item(508, 1).Value = "2024-12-05 15:47:17"
$ws.Cells.Item(508, 2).Value = "020920.NC"
$ws.Cells.Item(508, 3).Value = 596
$ws.Cells.Item(508, 4).Value = 183
$ws.Cells.Item(508, 5).Value = 32
$ws.Cells.Item(508, 6).Value = 1

$ws.Rows.Item(509).RowHeight = 13.55
$ws.Cells.Item(509, 1).Value = "2024-12-05 16:09:58"
$ws.Cells.Item(509, 2).Value = "020919.NC"
$ws.Cells.Item(509, 3).Value = 559
$ws.Cells.Item(509, 4).Value = 167
$ws.Cells.Item(509, 5).Value = 32
$ws.Cells.Item(509, 6).Value = 1

$ws.Rows.Item(510).RowHeight = 13.55
$ws.Cells.Item(510, 1).Value = "2024-12-06 08:06:34"
$ws.Cells.Item(510, 2).Value = "020918.NC"
$ws.Cells.Item(510, 3).Value = 398
$ws.Cells.Item(510, 4).Value = 195
$ws.Cells.Item(510, 5).Value = 32
$ws.Cells.Item(510, 6).Value = 1

$ws.Rows.Item(511).RowHeight = 13.55
$ws.Cells.Item(511, 1).Value = "2024-12-09 16:36:33"
$ws.Cells.Item(511, 2).Value = "020921.NC"
$ws.Cells.Item(511, 3).Value = 451
$ws.Cells.Item(511, 4).Value = 170
$ws.Cells.Item(511, 5).Value = 32
$ws.Cells.Item(511, 6).Value = 1

$ws.Rows.Item(512).RowHeight = 13.55
$ws.Cells.Item(512, 1).Value = "2024-12-10 11:10:49"
$ws.Cells.Item(512, 2).Value = "020911.NC"
$ws.Cells.Item(512, 3).Value = 531
$ws.Cells.Item(512, 4).Value = 155
$ws.Cells.Item(512, 5).Value = 32
$ws.Cells.Item(512, 6).Value = 1

$ws.Rows.Item(513).RowHeight = 13.55
$ws.Cells.Item(513, 1).Value = "2024-12-11 03:52:57"
$ws.Cells.Item(513, 2).Value = "020912.NC"
$ws.Cells.Item(513, 3).Value = 589
$ws.Cells.Item(513, 4).Value = 160
$ws.Cells.Item(513, 5).Value = 32
$ws.Cells.Item(513, 6).Value = 1

$ws.Rows.Item(514).RowHeight = 13.55
$ws.Cells.Item(514, 1).Value = "2024-12-12 00:45:31"
$ws.Cells.Item(514, 2).Value = "021907.NC"
$ws.Cells.Item(514, 3).Value = 8
$ws.Cells.Item(514, 4).Value = 15
$ws.Cells.Item(514, 5).Value = 10
$ws.Cells.Item(514, 6).Value = 1

$ws.Rows.Item(515).RowHeight = 13.55
$ws.Cells.Item(515, 1).Value = "2024-12-12 01:08:04"
$ws.Cells.Item(515, 2).Value = "021825.NC"
$ws.Cells.Item(515, 3).Value = 10
$ws.Cells.Item(515, 4).Value = 122
$ws.Cells.Item(515, 5).Value = 6
$ws.Cells.Item(515, 6).Value = 1

$ws.Rows.Item(516).RowHeight = 13.55
$ws.Cells.Item(516, 1).Value = "2024-12-12 02:23:21"
$ws.Cells.Item(516, 2).Value = "020914.NC"
$ws.Cells.Item(516, 3).Value = 576
$ws.Cells.Item(516, 4).Value = 163
$ws.Cells.Item(516, 5).Value = 32
$ws.Cells.Item(516, 6).Value = 1

$ws.Rows.Item(517).RowHeight = 13.55
$ws.Cells.Item(517, 1).Value = "2024-12-02 13:31:32"
$ws.Cells.Item(517, 2).Value = "020923.NC"
$ws.Cells.Item(517, 3).Value = 616
$ws.Cells.Item(517, 4).Value = 170
$ws.Cells.Item(517, 5).Value = 38
$ws.Cells.Item(517, 6).Value = 1

$ws.Rows.Item(518).RowHeight = 13.55
$ws.Cells.Item(518, 1).Value = "2024-12-03 08:46:46"
$ws.Cells.Item(518, 2).Value = "021532.NC"
$ws.Cells.Item(518, 3).Value = 12
$ws.Cells.Item(518, 4).Value = 27
$ws.Cells.Item(518, 5).Value = 12
$ws.Cells.Item(518, 6).Value = 1

$ws.Rows.Item(519).RowHeight = 13.55
$ws.Cells.Item(519, 1).Value = "2024-12-03 16:45:04"
$ws.Cells.Item(519, 2).Value = "020924.NC"
$ws.Cells.Item(519, 3).Value = 559
$ws.Cells.Item(519, 4).Value = 155
$ws.Cells.Item(519, 5).Value = 38
$ws.Cells.Item(519, 6).Value = 1

$ws.Rows.Item(520).RowHeight = 13.55
$ws.Cells.Item(520, 1).Value = "2024-12-04 08:57:31"
$ws.Cells.Item(520, 2).Value = "020925.NC"
$ws.Cells.Item(520, 3).Value = 314
$ws.Cells.Item(520, 4).Value = 87
$ws.Cells.Item(520, 5).Value = 38
$ws.Cells.Item(520, 6).Value = 1

Write-Output "done"